# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" stat. Update values for rows 2-10 per the new save_data regeneration.
$newValues = @{
    2  = 1
    3  = 5
    4  = 4
    5  = 3
    6  = 8
    7  = 3
    8  = 1
    9  = 4
    10 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
